{"js": "// Fix a typo: \"archetectual\" -> \"architectural\" in the Application\n// Structure section.\nconst typoResults = context.document.body.search(\"archetectual\", { matchCase: true });\ntypoResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < typoResults.items.length; i++) {\n  typoResults.items[i].insertText(\"architectural\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Update the \"Effort Spent\" section to describe the actual work split\n// instead of the generic 50/50 statement.\nconst effortResults = context.document.body.search(\n  \"we have both each done a 50% split of the work\",\n  { matchCase: true }\n);\neffortResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < effortResults.items.length; i++) {\n  effortResults.items[i].insertText(\n    \"Caleb completed 55% of the work and Quinn completed 45%\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Fix a typo: \"archetectual\" -> \"architectural\" in the Application\n# Structure section.\n$find = $d.Content.Find\n$find.Text = \"archetectual\"\n$find.Replacement.Text = \"architectural\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# Update the \"Effort Spent\" section to describe the actual work split\n# instead of the generic 50/50 statement.\n$find2 = $d.Content.Find\n$find2.Text = \"we have both each done a 50% split of the work\"\n$find2.Replacement.Text = \"Caleb completed 55% of the work and Quinn completed 45%\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
